$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update CNN+LSTM / SGD block (rows 6-8) with results from re-running on cluster (Tesla T4)
$ws.Range("A6").RowHeight = 14.4
$ws.Range("D6").Value = "210min"
$ws.Range("F6").Value = 81.900000000000006
$ws.Range("G6").Value = 0.55000000000000004
$ws.Range("H6").Value = "Run on Tesla T4 24GB"

$ws.Range("D7").Value = "210s per epoch"
$ws.Range("F7").Value = 69.8
$ws.Range("G7").Value = 0.9

$ws.Range("D8").Value = "52.1s"
$ws.Range("F8").Value = 71.2

# Update selection to match final saved state
$ws.Range("I12").Select()
